$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86:C86").NumberFormat = "@"

$ws.Range("A86").Value = "2025-10-18"
$ws.Range("B86").Value = "21:21:19"
$ws.Range("C86").Value = "1.00 EUR = 1,703.9130"
